$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC31_Verify_login")

$ws.Rows.Item(4).Insert()

$ws.Range("B4").Value = "CLICK"
$ws.Range("C4").Value = "LoginURL"
$ws.Range("D4").Value = "CSS"

$ws.Range("C5").Value = "Uname"
$ws.Range("C6").Value = "Password"
$ws.Range("C7").Value = "LoginButton"
